$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price-report row is inserted at row 82 ("Hortaliza, Vega
# Monumental Concepción - Papa"); every existing data row from 82..112
# shifts down by one (to 83..113), and the sheet's used range grows from
# A1:R112 to A1:R113.
$ws.Rows("82:82").Insert()

# Populate the newly inserted row 82 with the new week's record.
$ws.Range("A82").Value2 = 11
$ws.Range("B82").Value2 = "Vega Monumental Concepción"
$ws.Range("C82").Value2 = "Bíobío"
$ws.Range("D82").Value2 = 44460
$ws.Range("E82").Value2 = 8
$ws.Range("F82").Value2 = 100114001
$ws.Range("G82").Value2 = "Papa"
$ws.Range("H82").Value2 = "Asterix"
$ws.Range("I82").Value2 = "1a (guarda)"
$ws.Range("J82").Value2 = 2000
$ws.Range("K82").Value2 = 9500
$ws.Range("L82").Value2 = 10000
$ws.Range("M82").Value2 = 9750
$ws.Range("N82").Value2 = '$/saco 25 kilos'
$ws.Range("O82").Value2 = "Provincia de Arauco"
$ws.Range("P82").Value2 = 390
$ws.Range("Q82").Value2 = 25
$ws.Range("R82").Value2 = "Hortaliza"
